# Updated cryptos list - apply price/volume/coin changes per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.783.99'
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").Value = '2.289.59'
$ws.Range("E3").Value = '  -1.12%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '103.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '270.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("E7").Value = '  -0.40%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.606'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.95'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0933'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.44%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.95'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.28%  '

$ws.Range("E13").Value = '  +1.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.58%  '

$ws.Range("E15").Value = '  -2.01%  '

$ws.Range("D16").Value = '2.288.27'
$ws.Range("E16").Value = '  -1.39%  '

$ws.Range("D17").Value = '43.678.38'

$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.53%  '

$ws.Range("E21").Value = '  +9.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '233.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +14.38%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.27%  '

$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("E26").Value = '  -1.50%  '

$ws.Range("E27").Value = '  -0.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '40.45'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.62%  '

$ws.Range("E29").Value = '  -2.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '177.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.79'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0900'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.75%  '

$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.97'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +12.92%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.46'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.27%  '

$ws.Range("E35").Value = '  -0.20%  '

$ws.Range("E36").Value = '  -0.24%  '

$ws.Range("E37").Value = '  -2.76%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.55'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.12%  '

$ws.Range("E40").Value = '  -0.79%  '

$ws.Range("E41").Value = '  +0.41%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.30'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.94%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.63%  '

$ws.Range("E46").Value = '  -1.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '99.03'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.47%  '

$ws.Range("E49").Value = '  +7.89%  '

$ws.Range("E50").Value = '  +10.70%  '

$ws.Range("D51").Value = '2.512.93'
$ws.Range("E51").Value = '  -1.27%  '
